$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (closest achievable via Excel pixel-quantized ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 14.7
$ws.Columns.Item(2).ColumnWidth = 14.7

# Updated cell values
$ws.Range("A1").Value = -0.24002744690420741
$ws.Range("B1").Value = 0.23966817983475153
$ws.Range("A2").Value = -0.13930244796758728
$ws.Range("B2").Value = 0.13848884348537993
$ws.Range("A3").Value = -0.035552154918033096
$ws.Range("B3").Value = 0.035500765946057555
$ws.Range("A4").Value = -0.097489253776263496
$ws.Range("B4").Value = 0.097248037770629736
$ws.Range("A5").Value = -0.091248038985110469
$ws.Range("B5").Value = 0.090786688870180932
$ws.Range("A6").Value = -0.060648447275621908
$ws.Range("B6").Value = 0.060603359376520949
$ws.Range("A7").Value = -0.040603360813589617
$ws.Range("B7").Value = 0.040529987100335063
$ws.Range("A8").Value = -0.020529988544763178
$ws.Range("B8").Value = 0.020493761643324504
$ws.Range("A9").Value = -0.014493762896138129
$ws.Range("B9").Value = 0.014467862636347917
$ws.Range("A10").Value = -0.0084678638931947603
$ws.Range("B10").Value = 0.0084697017466979219
$ws.Range("A11").Value = -0.0039697029825092045
$ws.Range("B11").Value = 0.003967394196383367
$ws.Range("A12").Value = 0.002032604546259531
$ws.Range("B12").Value = -0.0020792253076691303
$ws.Range("A13").Value = -0.052358499696292249
$ws.Range("B13").Value = 0.052251957895824042
$ws.Range("A14").Value = -0.040251959244563373
$ws.Range("B14").Value = 0.040165489787087338
$ws.Range("A15").Value = -0.021048940925786574
$ws.Range("B15").Value = 0.02102568113991321
$ws.Range("A16").Value = -0.015025682413853492
$ws.Range("B16").Value = 0.015003395654440421
$ws.Range("A17").Value = -0.0090033969332186103
$ws.Range("B17").Value = 0.0089999986779991659
$ws.Range("A18").Value = -0.0903657331072516
$ws.Range("B18").Value = 0.090288646181694077
$ws.Range("A19").Value = -0.081288647409323289
$ws.Range("B19").Value = 0.08066518127427269
$ws.Range("A20").Value = -0.071665182538044547
$ws.Range("B20").Value = 0.071532429433343481
$ws.Range("A21").Value = -0.062532430704665742
$ws.Range("B21").Value = 0.062351734190375208
$ws.Range("A22").Value = -0.093936689811840779
$ws.Range("B22").Value = 0.093628251216456349
$ws.Range("A23").Value = -0.084628252472212395
$ws.Range("B23").Value = 0.084125185485470766
$ws.Range("A24").Value = -0.04212518722744818
$ws.Range("B24").Value = 0.041999998248549275
$ws.Range("A25").Value = -0.024138242831330814
$ws.Range("B25").Value = 0.024138308909783746
$ws.Range("A26").Value = -0.018138310147097769
$ws.Range("B26").Value = 0.018137722870413597
$ws.Range("A27").Value = -0.012137724108513659
$ws.Range("B27").Value = 0.012124973627941849
$ws.Range("A28").Value = -0.0061249748694640616
$ws.Range("B28").Value = 0.0061226325129872805
$ws.Range("A29").Value = 0.0058773661597779636
$ws.Range("B29").Value = -0.005877008075522383
$ws.Range("A30").Value = 0.025877006635795397
$ws.Range("B30").Value = -0.026163820058247467
$ws.Range("A31").Value = -0.036330586049281521
$ws.Range("B31").Value = 0.036264868753528745
$ws.Range("A32").Value = -0.0060004600025749255
$ws.Range("B32").Value = 0.0059999987619150374
